$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.934.43"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "'3.519.15"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D5").Value = "'601.90"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'181.62"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'3.517.24"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.595"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D13").Value = "'4.127.40"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'32.21"
$ws.Range("E14").Value = "  +9.93%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'67.884.35"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'3.529.60"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "'14.47"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").Value = "'401.57"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "'8.01"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'74.03"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "'5.73"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").Value = "'10.51"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'6.31"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "'23.99"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").Value = "'163.20"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "'2.82"
$ws.Range("E41").Value = "  +8.44%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'2.903.32"
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("D45").Value = "'26.47"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "'0.0736"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").Value = "'27.01"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'42.38"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'351.83"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  -0.95%  "
